$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 611.4
$ws.Range("I28").Value = 558.75
$ws.Range("K28").Value = 558.75
$ws.Range("M28").Value = -73.75
$ws.Range("H100").Value = 862.8570999999999
$ws.Range("I100").Value = 862.8570999999999
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 862.8570999999999
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -321.8570999999999
$ws.Range("N100").ClearContents()
$ws.Range("H125").Value = 22500
$ws.Range("J125").Value = 22500
$ws.Range("L125").Value = 202500
$ws.Range("N125").Value = -207420
$ws.Range("H132").Value = 2607.3
$ws.Range("I132").Value = 2563.6667
$ws.Range("K132").Value = 7691.000100000001
$ws.Range("M132").Value = -5161.000100000001
$ws.Range("H137").Value = 4933
$ws.Range("I137").Value = 2733.1667
$ws.Range("J137").Value = 9332.666999999999
$ws.Range("K137").Value = 8199.500100000001
$ws.Range("L137").Value = 27998.001
$ws.Range("M137").Value = -5649.500100000001
$ws.Range("N137").Value = -33098.001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4281290.5
$ws.Range("I32").Value = 4121366.2
$ws.Range("K32").Value = 4121366.2
$ws.Range("M32").Value = -4121079.2
$ws.Range("H97").Value = 1049.6666
$ws.Range("I97").Value = 995.25
$ws.Range("J97").Value = 1485
$ws.Range("K97").Value = 995.25
$ws.Range("L97").Value = 1485
$ws.Range("M97").Value = -499.25
$ws.Range("N97").Value = -2477

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H104").Value = 40000
$ws.Range("J104").Value = 40000
$ws.Range("L104").Value = 40000
$ws.Range("N104").Value = -46988
$ws.Range("H105").Value = 2220
$ws.Range("I105").Value = 2175
$ws.Range("K105").Value = 2175
$ws.Range("M105").Value = -428
$ws.Range("H107").Value = 1398.8334
$ws.Range("I107").Value = 1328.4286
$ws.Range("J107").Value = 1497.4
$ws.Range("K107").Value = 1328.4286
$ws.Range("L107").Value = 1497.4
$ws.Range("M107").Value = 591.5714
$ws.Range("N107").Value = -5337.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2275.75
$ws.Range("I58").Value = 1936
$ws.Range("J58").Value = 2615.5
$ws.Range("K58").Value = 1936
$ws.Range("L58").Value = 2615.5
$ws.Range("M58").Value = -1733
$ws.Range("N58").Value = -3021.5
$ws.Range("H99").Value = 1982.25
$ws.Range("I99").Value = 1312
$ws.Range("J99").Value = 2078
$ws.Range("K99").Value = 1312
$ws.Range("L99").Value = 2078
$ws.Range("M99").Value = 186
$ws.Range("N99").Value = -5074
$ws.Range("H105").Value = 3651.7
$ws.Range("I105").Value = 2335.3333
$ws.Range("K105").Value = 2335.3333
$ws.Range("M105").Value = -588.3332999999998
$ws.Range("H122").Value = 1699.5
$ws.Range("I122").Value = 2050
$ws.Range("K122").Value = 6150
$ws.Range("M122").Value = -3700
$ws.Range("H123").Value = 44250
$ws.Range("J123").Value = 44250
$ws.Range("L123").Value = 44250
$ws.Range("N123").Value = -54050
$ws.Range("H126").Value = 1982.25
$ws.Range("I126").Value = 1312
$ws.Range("J126").Value = 2078
$ws.Range("K126").Value = 3936
$ws.Range("L126").Value = 6234
$ws.Range("M126").Value = -1466
$ws.Range("N126").Value = -11174
$ws.Range("H136").Value = 2275.75
$ws.Range("I136").Value = 1936
$ws.Range("J136").Value = 2615.5
$ws.Range("K136").Value = 5808
$ws.Range("L136").Value = 7846.5
$ws.Range("M136").Value = -3258
$ws.Range("N136").Value = -12946.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2320
$ws.Range("I81").Value = 2000
$ws.Range("J81").Value = 2400
$ws.Range("K81").Value = 6000
$ws.Range("L81").Value = 7200
$ws.Range("M81").Value = -4877
$ws.Range("N81").Value = -9446
$ws.Range("H84").Value = 2320
$ws.Range("I84").Value = 2000
$ws.Range("J84").Value = 2400
$ws.Range("K84").Value = 18000
$ws.Range("L84").Value = 21600
$ws.Range("M84").Value = -12384
$ws.Range("N84").Value = -32832
$ws.Range("H96").Value = 15000
$ws.Range("J96").Value = 15000
$ws.Range("L96").Value = 45000
$ws.Range("N96").Value = -49118
$ws.Range("H122").Value = 886.75
$ws.Range("J122").Value = 889
$ws.Range("L122").Value = 8001
$ws.Range("N122").Value = -12901

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1313044.2
$ws.Range("I3").Value = 1100510
$ws.Range("J3").Value = 1667268
$ws.Range("K3").Value = 1100510
$ws.Range("L3").Value = 1667268
$ws.Range("M3").Value = -1100394
$ws.Range("N3").Value = -1667500
$ws.Range("H12").Value = 15000
$ws.Range("J12").Value = 15000
$ws.Range("L12").Value = 15000
$ws.Range("N12").Value = -15280
$ws.Range("H80").Value = 1696.625
$ws.Range("I80").Value = 1661.3334
$ws.Range("J80").Value = 1802.5
$ws.Range("K80").Value = 1661.3334
$ws.Range("L80").Value = 1802.5
$ws.Range("M80").Value = -663.3334
$ws.Range("N80").Value = -3798.5
$ws.Range("H83").Value = 1696.625
$ws.Range("I83").Value = 1661.3334
$ws.Range("J83").Value = 1802.5
$ws.Range("K83").Value = 8306.666999999999
$ws.Range("L83").Value = 9012.5
$ws.Range("M83").Value = -3314.666999999999
$ws.Range("N83").Value = -18996.5
$ws.Range("H87").Value = 47700
$ws.Range("J87").Value = 47700
$ws.Range("L87").Value = 47700
$ws.Range("N87").Value = -50196
$ws.Range("H90").Value = 47700
$ws.Range("J90").Value = 47700
$ws.Range("L90").Value = 143100
$ws.Range("N90").Value = -155580
$ws.Range("H102").Value = 2104.6667
$ws.Range("I102").Value = 2240
$ws.Range("K102").Value = 2240
$ws.Range("M102").Value = -618
$ws.Range("H122").Value = 3440.0715
$ws.Range("I122").Value = 3596.75
$ws.Range("K122").Value = 10790.25
$ws.Range("M122").Value = -8340.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1388.5714
$ws.Range("I16").Value = 1388.5714
$ws.Range("K16").Value = 1388.5714
$ws.Range("M16").Value = -1218.5714
$ws.Range("H40").Value = 5000
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H122").Value = 8896.272000000001
$ws.Range("I122").Value = 10666.625
$ws.Range("K122").Value = 31999.875
$ws.Range("M122").Value = -29549.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 17682.666
$ws.Range("J45").Value = 17682.666
$ws.Range("L45").Value = 17682.666
$ws.Range("N45").Value = -18664.666
$ws.Range("H100").Value = 9091165
$ws.Range("I100").Value = 11111375
$ws.Range("K100").Value = 22222750
$ws.Range("M100").Value = -22222209
$ws.Range("H107").Value = 3533.3333
$ws.Range("J107").Value = 3400
$ws.Range("L107").Value = 10200
$ws.Range("N107").Value = -14040
$ws.Range("H122").Value = 10000
$ws.Range("J122").Value = 10000
$ws.Range("L122").Value = 30000
$ws.Range("N122").Value = -34900
